# 1. add crash column for result
# Populate the 3 new convertible-bond rows (entries 3-6) on sheet "23Q3IN"
# (rows 7-14 were previously blank placeholder rows).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("23Q3IN")
$ws.Activate()

# --- Row 7/8 : 海波转债 ---------------------------------------------------
$ws.Range("A7").Value = "海波转债"
$ws.Range("B7").Value = "张海波"
$ws.Range("C7").Value = "桥梁设计"
$ws.Range("E7").Value = "一带一路"
$ws.Range("F7").Value = 3.16
$ws.Range("H7").Value = "未下修"
$ws.Range("J7").Value = -0.05
$ws.Range("K7").Value = "0.7/0.062"
$ws.Range("M7").Formula = "=SUM(B8:L8)"

$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 3
$ws.Range("H8").Value = 1

# --- Row 9/10 : 乐歌转债 --------------------------------------------------
$ws.Range("A9").Value = "乐歌转债"
$ws.Range("B9").Value = "项乐宏"
$ws.Range("C9").Value = "家用轻工"
$ws.Range("D9").Value = "明汯投资"
$ws.Range("E9").Value = "智能家居"
$ws.Range("F9").Value = 2.7
$ws.Range("H9").Value = "未下修"
$ws.Range("J9").Value = 0.27
$ws.Range("K9").Value = "15.16/5.17"
$ws.Range("M9").Formula = "=SUM(B10:L10)"

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 2
$ws.Range("H10").Value = 1

# --- Row 11/12 : 法兰转债 -------------------------------------------------
$ws.Range("A11").Value = "法兰转债"
$ws.Range("B11").Value = "质押担保/陶峰华"
$ws.Range("C11").Value = "起重机"
$ws.Range("E11").Value = "换电"
$ws.Range("F11").Value = 2.8
$ws.Range("H11").Value = "未下修"
$ws.Range("J11").Value = 0.11
$ws.Range("K11").Value = "3.09/1.39"
$ws.Range("M11").Formula = "=SUM(B12:L12)"

$ws.Range("B12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("H12").Value = 1

# --- Row 13/14 : 银信转债 -------------------------------------------------
$ws.Range("A13").Value = "银信转债"
$ws.Range("B13").Value = "詹立雄"
$ws.Range("C13").Value = "软件"
$ws.Range("D13").Value = "李怡名"
$ws.Range("E13").Value = "数据中心"
$ws.Range("F13").Value = 2.5
$ws.Range("H13").Value = "未下修"
$ws.Range("J13").Value = 0.14
$ws.Range("K13").Value = "5.53/1.19"
$ws.Range("M13").Formula = "=SUM(B14:L14)"

$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 2
$ws.Range("H14").Value = 1

# Scroll/selection bookkeeping so the active window matches the authored
# commit (top row 3, active cell J14).
$aw = $excel.ActiveWindow
$aw.ScrollRow = 3
$aw.ScrollColumn = 1
$ws.Range("J14").Select()

# --- sheet "22Q3IN" view bookkeeping -------------------------------------
$ws2 = $wb.Worksheets.Item("22Q3IN")
$ws2.Activate()
$aw2 = $excel.ActiveWindow
$aw2.ScrollRow = 16
$aw2.ScrollColumn = 1
$ws2.Range("D21").Select()

# Leave "23Q3IN" as the active/selected sheet, matching the source file.
$ws.Activate()
